# Update for the Documents 23.01.18
# - Set the Product Burndown entry for week of row 46 (column E, "Hinzugefuegt"/added
#   story points) to 18. This ripples through the running-total column I (I47:I55)
#   and the overall average in E56 via existing formulas, so no other cells need to
#   be written explicitly.
# - Move the active selection from E46 to E47 to match where editing continued.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E46").Value = 18

$ws.Range("E47").Select()
